$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 2).Value = 'New Zealand, Australia'
$ws.Cells.Item(7, 2).Value = 'New Zealand, Australia'
$ws.Cells.Item(13, 2).Value = 'Hong Kong, South Korea, Taiwan, Singapore'
$ws.Cells.Item(14, 2).Value = 'Palestine, North Africa, Algeria, Bahrain, Libya, Morocco, Arabia, Qatar, Somalia, Djibouti, Mauritania, United Arab Emirates, Sudan, Lebanon, Iraq, Comoros, Oman, Tunisia, Syria, Kuwait, Saudi Arabia, Jordan, Yemen, Egypt'
$ws.Cells.Item(17, 2).Value = 'South Korea, Japan, New Zealand, India, China, Australia'
$ws.Cells.Item(20, 2).Value = 'Lithuania, Latvia, Estonia'
$ws.Cells.Item(21, 2).Value = 'Brazil, India, South Africa, China'
$ws.Cells.Item(22, 2).Value = 'Luxembourg, Belgium'
$ws.Cells.Item(23, 2).Value = 'Italy, Germany, France, Spain'
$ws.Cells.Item(24, 2).Value = 'India, Myanmar, Sri Lanka, Bhutan, Thailand, Nepal, Bangladesh'
$ws.Cells.Item(25, 2).Value = 'Brazil, India, Russia, China'
$ws.Cells.Item(26, 2).Value = 'India, China, Brazil, Russia, South Africa'
$ws.Cells.Item(30, 2).Value = 'New Zealand, Canada, Australia'
$ws.Cells.Item(32, 2).Value = 'Serbia, Montenegro, Kosovo, Bosnia-Herzegovina, North Macedonia, Albania, Moldova'
$ws.Cells.Item(34, 2).Value = 'Colombia, Vietnam, Indonesia, Egypt, Turkey, South Africa'
$ws.Cells.Item(35, 2).Value = 'Vietnam, Myanmar, Cambodia, Laos'
$ws.Cells.Item(36, 2).Value = 'Hungary, Vietnam, Czechoslovakia, Cuba, Romania, Bulgaria, East Germany, Poland, Mongolia'
$ws.Cells.Item(38, 2).Value = 'Kyrgyzstan, Kazakhstan, Uzbekistan, Republics, Armenia, Azerbaijan, Belarus, Tajikistan, Russia, Turkmenistan, Moldova'
$ws.Cells.Item(41, 2).Value = 'Kyrgyzstan, Serbia, Kazakhstan, Armenia, Belarus, Tajikistan, Russia, Afghanistan'
$ws.Cells.Item(45, 2).Value = 'New Zealand, Slovakia, Denmark, Finland, Ireland, Portugal, Luxembourg, Japan, Belgium, Czech Republic, Slovenia, Poland, Italy, Canada, Sweden, Switzerland, Hungary, Netherlands, Austria, Iceland, Australia, United States, South Korea, France, Norway, Spain, Greece, United Kingdom, Germany'
$ws.Cells.Item(46, 2).Value = 'Austria, Germany, Switzerland'
$ws.Cells.Item(47, 2).Value = 'Austria, Germany, Switzerland'
$ws.Cells.Item(48, 2).Value = 'South Sudan, Tanzania, Burundi, Uganda, Kenya, Rwanda'
$ws.Cells.Item(49, 2).Value = 'Ukraine, Armenia, Azerbaijan, Belarus, Moldova, Georgia'
$ws.Cells.Item(50, 2).Value = 'Kyrgyzstan, Turkmenistan, Kazakhstan, Uzbekistan, Iran, Azerbaijan, Turkey, Tajikistan, Afghanistan, Pakistan'
$ws.Cells.Item(52, 2).Value = 'Kyrgyzstan, Cuba, Kazakhstan, Uzbekistan, Armenia, Belarus, Russia, Moldova'
$ws.Cells.Item(60, 2).Value = 'New Zealand, Canada, Australia'
$ws.Cells.Item(62, 2).Value = 'Hong Kong, South Korea, Taiwan, Singapore'
$ws.Cells.Item(66, 2).Value = 'Italy, France'
$ws.Cells.Item(67, 2).Value = 'Italy, France, Spain'
$ws.Cells.Item(68, 2).Value = 'Brazil, India, Germany, Japan'
$ws.Cells.Item(71, 2).Value = 'Poland, Italy, Germany, Spain'
$ws.Cells.Item(72, 2).Value = 'France, Japan, Italy, Germany, Canada'
$ws.Cells.Item(73, 2).Value = 'India, China, Brazil, Mexico, South Africa'
$ws.Cells.Item(74, 2).Value = 'India, Argentina, Japan, Indonesia, Turkey, Mexico, Italy, Canada, Australia, Brazil, United States, Russia, South Korea, Saudi Arabia, France, United Kingdom, China, South Africa, Germany'
$ws.Cells.Item(76, 2).Value = 'Moldova, Ukraine, Azerbaijan, Georgia'
$ws.Cells.Item(77, 2).Value = 'Kuwait, Iraq, Saudi Arabia, Bahrain, Oman, Qatar, UAE'
$ws.Cells.Item(78, 2).Value = 'Hong Kong, China, Taiwan, Macau'
$ws.Cells.Item(79, 2).Value = 'Austria, Germany, Switzerland'
$ws.Cells.Item(80, 2).Value = 'Brazil, India, South Africa'
$ws.Cells.Item(88, 2).Value = 'Palestine, Israel, Jordan, Lebanon, Syria'
$ws.Cells.Item(89, 2).Value = 'Brazil, Argentina, Uruguay, Paraguay'
$ws.Cells.Item(90, 2).Value = 'Indonesia, Republic of Korea, Turkey, Australia, Mexico'
$ws.Cells.Item(91, 2).Value = 'Mexico, Turkey, Nigeria, Indonesia'
$ws.Cells.Item(102, 2).Value = 'Mexico, Canada'
$ws.Cells.Item(103, 2).Value = 'South Korea, Vietnam, Indonesia, Nigeria, Egypt, Philippines, Iran, Turkey, Bangladesh, Mexico, Pakistan'
$ws.Cells.Item(106, 2).Value = 'Mexico, United States, Canada'
$ws.Cells.Item(107, 2).Value = 'Sweden, Norway, Denmark, Finland, Iceland'
$ws.Cells.Item(115, 2).Value = 'Russia, China, France'
$ws.Cells.Item(116, 2).Value = 'Colombia, Mexico, Chile, Peru'
$ws.Cells.Item(117, 2).Value = 'Colombia, Mexico, Chile, Peru'
$ws.Cells.Item(119, 2).Value = 'Ireland, Portugal, Spain, Greece, Italy'
$ws.Cells.Item(124, 2).Value = 'Norway, Sweden, Denmark, Finland'
$ws.Cells.Item(125, 2).Value = 'Kyrgyzstan, India, Kazakhstan, Uzbekistan, China, Tajikistan, Russia, Pakistan'
$ws.Cells.Item(126, 2).Value = 'Luxembourg, Lorraine, Saarland'
$ws.Cells.Item(128, 2).Value = 'Argentina, Uruguay, Chile, Paraguay'
$ws.Cells.Item(129, 2).Value = 'India, Sri Lanka, Bhutan, Maldives, Nepal, Bangladesh, Afghanistan, Pakistan'
$ws.Cells.Item(130, 2).Value = 'Vietnam, Indonesia, Philippines, Thailand, Malaysia'
$ws.Cells.Item(131, 2).Value = 'Kyrgyzstan, Hungary, Kazakhstan, Uzbekistan, Azerbaijan, Turkey'
$ws.Cells.Item(132, 2).Value = 'Wales, UK, Northern Ireland, United Kingdom, Scotland, England'
$ws.Cells.Item(135, 2).Value = 'Poland, Slovakia, Hungary'
$ws.Cells.Item(136, 2).Value = 'Vietnam, Indonesia, Argentina, Turkey, South Africa'
